$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "hhh"
$ws.Range("B5").Value = "asd"
$ws.Range("C5").Value = "Monumento"
$ws.Range("D5").Value = "24 horas"
$ws.Range("E5").Value = 5
$ws.Range("F5").Value = "dxhdsf"
$ws.Range("G5").Value = "34234"
$ws.Range("H5").Value = "sadgad"

$ws.Range("A6").Value = "jjj"
$ws.Range("B6").Value = "uiy"
$ws.Range("C6").Value = "Parque"
$ws.Range("D6").Value = "24 horas"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = "ey536"
$ws.Range("G6").Value = "454564"
$ws.Range("H6").Value = "erytewt"
